# Apply updated cryptocurrency price/volume data to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.325.99'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.931.39'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.05'
$ws.Range('E5').Value = '  +2.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.7143'
$ws.Range('E6').Value = '  -0.84%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '27.47'
$ws.Range('E9').Value = '  +4.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07191'
$ws.Range('E10').Value = '  +5.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7989'
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08085'
$ws.Range('E12').Value = '  +1.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.929.17'
$ws.Range('E13').Value = '  +0.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.427'
$ws.Range('E14').Value = '  +0.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.79'
$ws.Range('E15').Value = '  +0.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.84'
$ws.Range('E16').Value = '  +2.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.293.30'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '251.24'
$ws.Range('E18').Value = '  -3.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008131'
$ws.Range('E19').Value = '  +2.39%  '
$ws.Range('E20').Value = '  -0.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.178.86'
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9995'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.912'
$ws.Range('E24').Value = '  +0.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.722'
$ws.Range('E25').Value = '  +0.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.84'
$ws.Range('E26').Value = '  +3.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.22'
$ws.Range('E27').Value = '  +1.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.314'
$ws.Range('E28').Value = '  +1.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1284'
$ws.Range('E29').Value = '  -3.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.361'
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('E31').Value = '  -0.23%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.430'
$ws.Range('E32').Value = '  +0.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.208'
$ws.Range('E33').Value = '  +0.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05211'
$ws.Range('E34').Value = '  +2.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.264'
$ws.Range('E35').Value = '  +5.86%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7473'
$ws.Range('E36').Value = '  +1.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.762'
$ws.Range('E37').Value = '  +1.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01960'
$ws.Range('E38').Value = '  +1.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.797'
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '78.86'
$ws.Range('E40').Value = '  -1.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.424'
$ws.Range('E41').Value = '  -2.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4528'
$ws.Range('E42').Value = '  +1.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.025'
$ws.Range('E43').Value = '  +1.06%  '
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8401'
$ws.Range('E45').Value = '  +0.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.74'
$ws.Range('E46').Value = '  -0.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.780'
$ws.Range('E47').Value = '  +0.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.403'
$ws.Range('E48').Value = '  +1.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '36.61'
$ws.Range('E49').Value = '  +1.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06100'
$ws.Range('E50').Value = '  +3.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4175'
$ws.Range('E51').Value = '  +1.73%  '
